{"js": "// Replace the three-digit division answers with the new set of values.\n// Each old string is unique in the document, so a simple matchCase search\n// + full-text replace is safe and format-preserving (keeps run formatting).\nconst replacements = [\n  [\"802\u00f78=100, 2\", \"765\u00f77=109, 2\"],\n  [\"911\u00f76=151, 5\", \"634\u00f75=126, 4\"],\n  [\"717\u00f74=179, 1\", \"686\u00f77=98, 0\"],\n  [\"367\u00f76=61, 1\", \"317\u00f76=52, 5\"],\n  [\"241\u00f74=60, 1\", \"923\u00f75=184, 3\"],\n  [\"136\u00f78=17, 0\", \"543\u00f75=108, 3\"],\n  [\"829\u00f77=118, 3\", \"956\u00f76=159, 2\"],\n  [\"228\u00f78=28, 4\", \"973\u00f77=139, 0\"],\n  [\"204\u00f77=29, 1\", \"577\u00f72=288, 1\"],\n  [\"734\u00f74=183, 2\", \"182\u00f79=20, 2\"],\n  [\"585\u00f73=195, 0\", \"934\u00f78=116, 6\"],\n  [\"471\u00f73=157, 0\", \"919\u00f78=114, 7\"],\n  [\"475\u00f78=59, 3\", \"940\u00f77=134, 2\"],\n  [\"281\u00f75=56, 1\", \"678\u00f74=169, 2\"],\n  [\"701\u00f78=87, 5\", \"931\u00f73=310, 1\"],\n  [\"866\u00f75=173, 1\", \"219\u00f75=43, 4\"],\n  [\"192\u00f73=64, 0\", \"140\u00f78=17, 4\"],\n  [\"457\u00f78=57, 1\", \"694\u00f73=231, 1\"],\n  [\"855\u00f78=106, 7\", \"987\u00f77=141, 0\"],\n  [\"217\u00f75=43, 2\", \"901\u00f79=100, 1\"],\n  [\"801\u00f76=133, 3\", \"421\u00f72=210, 1\"],\n  [\"200\u00f74=50, 0\", \"877\u00f74=219, 1\"],\n  [\"472\u00f72=236, 0\", \"497\u00f75=99, 2\"],\n  [\"809\u00f73=269, 2\", \"933\u00f76=155, 3\"],\n  [\"995\u00f79=110, 5\", \"613\u00f72=306, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit division answers with the new set of values.\n# Each \"old\" string is unique in the document, so Find/Replace (MatchCase,\n# whole-document scope, ReplaceAll) is safe and keeps the existing run\n# formatting (font/size) untouched.\n$pairs = @(\n  @(\"802\u00f78=100, 2\", \"765\u00f77=109, 2\"),\n  @(\"911\u00f76=151, 5\", \"634\u00f75=126, 4\"),\n  @(\"717\u00f74=179, 1\", \"686\u00f77=98, 0\"),\n  @(\"367\u00f76=61, 1\", \"317\u00f76=52, 5\"),\n  @(\"241\u00f74=60, 1\", \"923\u00f75=184, 3\"),\n  @(\"136\u00f78=17, 0\", \"543\u00f75=108, 3\"),\n  @(\"829\u00f77=118, 3\", \"956\u00f76=159, 2\"),\n  @(\"228\u00f78=28, 4\", \"973\u00f77=139, 0\"),\n  @(\"204\u00f77=29, 1\", \"577\u00f72=288, 1\"),\n  @(\"734\u00f74=183, 2\", \"182\u00f79=20, 2\"),\n  @(\"585\u00f73=195, 0\", \"934\u00f78=116, 6\"),\n  @(\"471\u00f73=157, 0\", \"919\u00f78=114, 7\"),\n  @(\"475\u00f78=59, 3\", \"940\u00f77=134, 2\"),\n  @(\"281\u00f75=56, 1\", \"678\u00f74=169, 2\"),\n  @(\"701\u00f78=87, 5\", \"931\u00f73=310, 1\"),\n  @(\"866\u00f75=173, 1\", \"219\u00f75=43, 4\"),\n  @(\"192\u00f73=64, 0\", \"140\u00f78=17, 4\"),\n  @(\"457\u00f78=57, 1\", \"694\u00f73=231, 1\"),\n  @(\"855\u00f78=106, 7\", \"987\u00f77=141, 0\"),\n  @(\"217\u00f75=43, 2\", \"901\u00f79=100, 1\"),\n  @(\"801\u00f76=133, 3\", \"421\u00f72=210, 1\"),\n  @(\"200\u00f74=50, 0\", \"877\u00f74=219, 1\"),\n  @(\"472\u00f72=236, 0\", \"497\u00f75=99, 2\"),\n  @(\"809\u00f73=269, 2\", \"933\u00f76=155, 3\"),\n  @(\"995\u00f79=110, 5\", \"613\u00f72=306, 1\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute([ref]$null, [ref]$true, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]2)\n}\n"}
